$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order "Recorded By" email lists (same people, new order) ---
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G5").Value = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# --- Updated statistics (recorded/missing session counts & coverage %) ---
# Leading "'" keeps these as literal text (e.g. "48.3%") instead of letting
# Excel auto-convert the percent-looking string into a numeric percentage.
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 1
$ws.Range("L9").Value = "'48.3%"
$ws.Range("L10").Value = "'24.7%"

# Mirror the same summary numbers in the PARASITOLOGY roll-up row
$ws.Range("O15").Value = 14
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "'48.3%"
$ws.Range("S15").Value = "'24.7%"

# --- Row 22 (PATHOLOGY LAB/MUSEUM session) just got recorded ---
# Re-format A22:I22 to the "Recorded" (green) look used by other recorded rows
$ws.Range("A2:I2").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G22").Value = "Alaa.A.Mostafa@med.asu.edu.eg"
$ws.Range("H22").Value = "19/251"
$ws.Range("I22").Value = "Recorded"
